$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scales")

# Row 10 / Row 11 data used to start at column B (label) / C (first value).
# Migrate the whole "db -> amplitude" mini-table one column to the right
# (B->C, C->D, ..., O->P), turning the old formula in C11 into a label.

$ws.Range("B10").ClearContents()

$ws.Range("C10").Value = "db"
$ws.Range("C11").Value = "Amplitude"

$dbValues = @(0, -5, -10, -15, -20, -25, -30, -35, -40, -45, -50, -55, -60)
$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $ws.Range($col + "10").Value = $dbValues[$i]
    $ws.Range($col + "11").Formula = "=POWER(10, " + $col + "10/20)"
}

$ws.Range("L11").Select()
